# Refresh the cryptos price/volume table (and a couple of re-ranked rows)
# to match the latest scrape, cell by cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row lists only the columns that actually changed for that coin:
# B/C = name+link (only differ where rows got re-ranked), D = price, E = 1h change.
$updates = @(
    @{ Row = 2; Cells = @{ 'D'='61.043.90'; 'E'='  -1.81%  ' } }
    @{ Row = 3; Cells = @{ 'D'='3.386.04'; 'E'='  -0.83%  ' } }
    @{ Row = 4; Cells = @{ 'D'='0.999'; 'E'='  -0.09%  ' } }
    @{ Row = 5; Cells = @{ 'D'='574.29'; 'E'='  -1.64%  ' } }
    @{ Row = 6; Cells = @{ 'D'='136.79'; 'E'='  -1.62%  ' } }
    @{ Row = 7; Cells = @{ 'E'='  +0.03%  ' } }
    @{ Row = 8; Cells = @{ 'D'='3.383.03'; 'E'='  -0.86%  ' } }
    @{ Row = 9; Cells = @{ 'E'='  -1.53%  ' } }
    @{ Row = 10; Cells = @{ 'D'='7.68'; 'E'='  +1.55%  ' } }
    @{ Row = 11; Cells = @{ 'D'='0.123'; 'E'='  -4.06%  ' } }
    @{ Row = 12; Cells = @{ 'E'='  -2.82%  ' } }
    @{ Row = 13; Cells = @{ 'D'='3.959.87'; 'E'='  -0.92%  ' } }
    @{ Row = 14; Cells = @{ 'E'='  +0.87%  ' } }
    @{ Row = 15; Cells = @{ 'E'='  -4.33%  ' } }
    @{ Row = 16; Cells = @{ 'D'='3.383.85'; 'E'='  -0.84%  ' } }
    @{ Row = 17; Cells = @{ 'D'='25.46' } }
    @{ Row = 18; Cells = @{ 'D'='61.155.99'; 'E'='  -1.78%  ' } }
    @{ Row = 19; Cells = @{ 'D'='13.83'; 'E'='  -2.73%  ' } }
    @{ Row = 20; Cells = @{ 'D'='5.76'; 'E'='  -1.62%  ' } }
    @{ Row = 21; Cells = @{ 'D'='9.35'; 'E'='  -2.75%  ' } }
    @{ Row = 22; Cells = @{ 'D'='375.89'; 'E'='  -5.18%  ' } }
    @{ Row = 23; Cells = @{ 'D'='3.515.65'; 'E'='  -1.01%  ' } }
    @{ Row = 24; Cells = @{ 'D'='0.551'; 'E'='  -3.19%  ' } }
    @{ Row = 25; Cells = @{ 'E'='  +0.26%  ' } }
    @{ Row = 26; Cells = @{ 'E'='  -4.82%  ' } }
    @{ Row = 27; Cells = @{ 'D'='71.01'; 'E'='  -1.05%  ' } }
    @{ Row = 28; Cells = @{ 'D'='0.182'; 'E'='  +12.30%  ' } }
    @{ Row = 29; Cells = @{ 'D'='1.66'; 'E'='  -1.93%  ' } }
    @{ Row = 30; Cells = @{ 'E'='  -0.08%  ' } }
    @{ Row = 31; Cells = @{ 'D'='7.40'; 'E'='  -4.85%  ' } }
    @{ Row = 32; Cells = @{ 'D'='8.07'; 'E'='  -2.74%  ' } }
    @{ Row = 33; Cells = @{ 'D'='2.15'; 'E'='  -2.43%  ' } }
    @{ Row = 34; Cells = @{ 'E'='  -0.04%  ' } }
    @{ Row = 35; Cells = @{ 'D'='23.53'; 'E'='  -0.32%  ' } }
    @{ Row = 36; Cells = @{ 'D'='5.18'; 'E'='  -5.07%  ' } }
    @{ Row = 37; Cells = @{ 'E'='  -3.23%  ' } }
    @{ Row = 38; Cells = @{ 'D'='6.79'; 'E'='  -2.63%  ' } }
    @{ Row = 39; Cells = @{ 'D'='164.29'; 'E'='  -0.32%  ' } }
    @{ Row = 40; Cells = @{ 'D'='0.0755'; 'E'='  -4.90%  ' } }
    @{ Row = 41; Cells = @{ 'B'='EnergySwap'; 'C'='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; 'D'='25.47'; 'E'='  +1.88%  ' } }
    @{ Row = 42; Cells = @{ 'B'='FirstDigitalUSD'; 'C'='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; 'D'='0.999'; 'E'='  -0.16%  ' } }
    @{ Row = 43; Cells = @{ 'B'='Mantle'; 'C'='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; 'D'='0.774'; 'E'='  -1.95%  ' } }
    @{ Row = 44; Cells = @{ 'B'='OKB'; 'C'='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; 'D'='41.74'; 'E'='  -0.32%  ' } }
    @{ Row = 45; Cells = @{ 'D'='1.70'; 'E'='  -6.36%  ' } }
    @{ Row = 46; Cells = @{ 'D'='1.20'; 'E'='  -6.73%  ' } }
    @{ Row = 47; Cells = @{ 'E'='  -2.91%  ' } }
    @{ Row = 48; Cells = @{ 'D'='2.481.69'; 'E'='  +4.85%  ' } }
    @{ Row = 49; Cells = @{ 'D'='6.79'; 'E'='  -2.21%  ' } }
    @{ Row = 50; Cells = @{ 'D'='22.82'; 'E'='  -3.17%  ' } }
    @{ Row = 51; Cells = @{ 'D'='2.42'; 'E'='  +3.19%  ' } }
)

foreach ($update in $updates) {
    foreach ($col in $update.Cells.Keys) {
        $cellRef = "$col$($update.Row)"
        $value = $update.Cells[$col]
        $range = $ws.Range($cellRef)
        if ($value -match '^[+-]?\d+(\.\d+)?$') {
            # Plain decimal-looking text (e.g. "7.40", "1.20"): force Text so Excel
            # does not silently reparse it as a Number and drop trailing zeros,
            # then drop the formatting override again so the cell keeps its
            # original (default) style once the value is safely stored as text.
            $range.NumberFormat = "@"
            $range.Value = $value
            $range.ClearFormats()
        } else {
            $range.Value = $value
        }
    }
}
